$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 13
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2831
# row 17
$ws.Range("H17").Value = 879.0513
$ws.Range("J17").Value = 879.0513
$ws.Range("L17").Value = 2637.1539
$ws.Range("N17").Value = -2973.1539
# row 40
$ws.Range("H40").Value = 2698.7778
$ws.Range("I40").Value = 1999.6666
$ws.Range("J40").Value = 3048.3333
$ws.Range("K40").Value = 1999.6666
$ws.Range("L40").Value = 3048.3333
$ws.Range("M40").Value = -1824.6666
$ws.Range("N40").Value = -3398.3333
# row 70
$ws.Range("H70").Value = 1500
$ws.Range("I70").Value = 1750
$ws.Range("J70").Value = 1000
$ws.Range("K70").Value = 5250
$ws.Range("L70").Value = 3000
$ws.Range("M70").Value = -4980
$ws.Range("N70").Value = -3540
# row 73
$ws.Range("H73").Value = 1500
$ws.Range("I73").Value = 1750
$ws.Range("J73").Value = 1000
$ws.Range("K73").Value = 5250
$ws.Range("L73").Value = 3000
$ws.Range("M73").Value = -4314
$ws.Range("N73").Value = -4872
# row 137
$ws.Range("H137").Value = 1252.2106
$ws.Range("I137").Value = 1236.9286
$ws.Range("J137").Value = 1295
$ws.Range("K137").Value = 3710.7858
$ws.Range("L137").Value = 3885
$ws.Range("M137").Value = -1160.7858
$ws.Range("N137").Value = -8985
# row 138
$ws.Range("H138").Value = 1344.746
$ws.Range("J138").Value = 1782.6
$ws.Range("L138").Value = 5347.799999999999
$ws.Range("N138").Value = -15627.8

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 9252.666999999999
$ws.Range("I2").Value = 889.875
$ws.Range("K2").Value = 889.875
$ws.Range("M2").Value = -776.875
# row 32
$ws.Range("H32").Value = 4125.8823
$ws.Range("I32").Value = 4125.8823
$ws.Range("K32").Value = 4125.8823
$ws.Range("M32").Value = -3838.8823
# row 116
$ws.Range("H116").Value = 9252.666999999999
$ws.Range("I116").Value = 889.875
$ws.Range("K116").Value = 889.875
$ws.Range("M116").Value = 1404.125

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 9252.666999999999
$ws.Range("I3").Value = 889.875
$ws.Range("K3").Value = 889.875
$ws.Range("M3").Value = -775.875
# row 5
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -87
$ws.Range("N5").Value = ""
# row 7
$ws.Range("H7").Value = 89
$ws.Range("I7").Value = 83.5
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 83.5
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 29.5
$ws.Range("N7").Value = -326
# row 80
$ws.Range("H80").Value = 855.4167
$ws.Range("I80").Value = 456
$ws.Range("J80").Value = 988.55554
$ws.Range("K80").Value = 456
$ws.Range("L80").Value = 988.55554
$ws.Range("M80").Value = 542
$ws.Range("N80").Value = -2984.55554
# row 83
$ws.Range("H83").Value = 855.4167
$ws.Range("I83").Value = 456
$ws.Range("J83").Value = 988.55554
$ws.Range("K83").Value = 2280
$ws.Range("L83").Value = 4942.7777
$ws.Range("M83").Value = 2712
$ws.Range("N83").Value = -14926.7777
# row 140
$ws.Range("H140").Value = 20780
$ws.Range("J140").Value = 20780
$ws.Range("L140").Value = 20780
$ws.Range("N140").Value = -31140

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1733.425
$ws.Range("I31").Value = 889.1667
$ws.Range("K31").Value = 889.1667
$ws.Range("M31").Value = -594.1667
# row 34
$ws.Range("H34").Value = 1733.425
$ws.Range("I34").Value = 889.1667
$ws.Range("K34").Value = 889.1667
$ws.Range("M34").Value = -687.1667
# row 99
$ws.Range("H99").Value = 1639.9565
$ws.Range("I99").Value = 1665.25
$ws.Range("J99").Value = 1612.3636
$ws.Range("K99").Value = 1665.25
$ws.Range("L99").Value = 1612.3636
$ws.Range("M99").Value = -167.25
$ws.Range("N99").Value = -4608.3636
# row 126
$ws.Range("H126").Value = 1639.9565
$ws.Range("I126").Value = 1665.25
$ws.Range("J126").Value = 1612.3636
$ws.Range("K126").Value = 4995.75
$ws.Range("L126").Value = 4837.0908
$ws.Range("M126").Value = -2525.75
$ws.Range("N126").Value = -9777.0908
# row 132
$ws.Range("H132").Value = 2723.0667
$ws.Range("I132").Value = 2139.4546
$ws.Range("J132").Value = 4328
$ws.Range("K132").Value = 6418.3638
$ws.Range("L132").Value = 12984
$ws.Range("M132").Value = -3888.3638
$ws.Range("N132").Value = -18044

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 1906.5555
$ws.Range("I132").Value = 1520
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4560
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2030
$ws.Range("N132").Value = -20057

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2181.125
$ws.Range("I7").Value = 2144.4
$ws.Range("K7").Value = 2144.4
$ws.Range("M7").Value = -2032.4
# row 55
$ws.Range("H55").Value = 322.625
$ws.Range("I55").Value = 91.5
$ws.Range("J55").Value = 1016
$ws.Range("K55").Value = 91.5
$ws.Range("L55").Value = 1016
$ws.Range("M55").Value = 81.5
$ws.Range("N55").Value = -1362
# row 61
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = ""
# row 113
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
# row 126
$ws.Range("H126").Value = 2181.125
$ws.Range("I126").Value = 2144.4
$ws.Range("K126").Value = 6433.200000000001
$ws.Range("M126").Value = -3963.200000000001
# row 132
$ws.Range("H132").Value = 28439.73
$ws.Range("I132").Value = 998.6429000000001
$ws.Range("J132").Value = 113812
$ws.Range("K132").Value = 2995.9287
$ws.Range("L132").Value = 341436
$ws.Range("M132").Value = -465.9287000000004
$ws.Range("N132").Value = -346496

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 64
$ws.Range("H64").Value = 18000
$ws.Range("J64").Value = 18000
$ws.Range("L64").Value = 18000
$ws.Range("N64").Value = -18496
# row 67
$ws.Range("H67").Value = 18000
$ws.Range("J67").Value = 18000
$ws.Range("L67").Value = 18000
$ws.Range("N67").Value = -19716
# row 132
$ws.Range("H132").Value = 3226.8462
$ws.Range("I132").Value = 2904.762
$ws.Range("J132").Value = 4579.6
$ws.Range("K132").Value = 8714.286
$ws.Range("L132").Value = 13738.8
$ws.Range("M132").Value = -6184.286
$ws.Range("N132").Value = -18798.8
